$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at the top of the "Femacal de La Calera - Limón" weekly block
# (before old row 761), shifting all existing data rows down by 4.
$ws.Rows("761:764").Insert()

# New data for the inserted rows (week of 44610), same market/product metadata
# as the rest of the block.
$newRows = @(
    @{ D=44610; L="1a amarillo"; M=56;  N=14000; O=14000; P=14000; S=875; T=16 },
    @{ D=44610; L="1a plateado"; M=153; N=14000; O=15000; P=14490; S=906; T=16 },
    @{ D=44610; L="2a amarillo"; M=50;  N=12000; O=12000; P=12000; S=750; T=16 },
    @{ D=44610; L="2a plateado"; M=167; N=12000; O=13000; P=12479; S=780; T=16 }
)

$r = 761
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = 3
    $ws.Cells.Item($r, 2).Value = "Femacal de La Calera"
    $ws.Cells.Item($r, 3).Value = "Coquimbo"
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = 5
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100102
    $ws.Cells.Item($r, 8).Value = "Cítricos"
    $ws.Cells.Item($r, 9).Value = 100102003
    $ws.Cells.Item($r, 10).Value = "Limón"
    $ws.Cells.Item($r, 11).Value = "Sin especificar"
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = "$/malla 16 kilos"
    $ws.Cells.Item($r, 18).Value = "Provincia de Quillota"
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
    $r = $r + 1
}
